$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and the row-50/51 coin swap)
# Columns D (Price) and E (Volume 1h) hold number-like text that Excel would
# otherwise auto-convert, so each touched cell is forced to Text format first.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.370.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.487.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.86"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.85"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.48%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.82%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.648"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.98"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.08%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.53"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.038.58"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "600.03"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.462.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.83"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.55"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.500.22"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.96%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.87%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.64"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +11.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.63"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.74%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.38%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.45%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.62"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.12"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +20.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.39"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.24"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -7.60%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "519.89"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.67"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.617.41"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.66"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.46%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.11%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.94"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.58%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.92%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.78"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.93%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "OceanProtocol"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.36"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -11.79%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000242"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -8.33%  "
